$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 (Hydrogen / Non-metallic minerals) no longer carries a value - clear it
# back to a blank cell, like the other untouched cells in the sheet.
$ws.Range("D3").Value = ""

# The old row 7 ("Other") is split in two: its total moves down to a new
# row 8 that keeps the "Other" label, while row 7 is repurposed for the
# newly reported "Biogas" category (currently zero).
$ws.Range("A7:D7").Copy($ws.Range("A8:D8")) | Out-Null

# Row 7 becomes "Biogas" with no demand recorded yet.
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 0

# Row 8 keeps the "Other" label (copied formatting/blank cells already in
# place) together with the total that used to sit on row 7.
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 576.0680523255862
